$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.47"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.56%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.705"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.31%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06176"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.45%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.60%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8512"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.60%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9118"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.92%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1401"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.38%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05151"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.37%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07097"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.00%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03111"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.14%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09046"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.14%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001537"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.45%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006168"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.64%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005982"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.11%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.453"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.14%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.173"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.80%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.14%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1310"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.35%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.096"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.23%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.46%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001179"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.51%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004047"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.56%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.06%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.13%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03986"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.11%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.01%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004139"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.37%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.32%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01323"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-18.45%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.17%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.06%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.02122"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-61.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2578"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "90.24%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.06%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
